$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1439688715953307
$ws.Range("C2").Value = 0.6381322957198443
$ws.Range("J2").Value = 0.003891050583657588
$ws.Range("P2").Value = 0.132295719844358
$ws.Range("S2").Value = 0.08171206225680934
$ws.Range("B3").Value = 0.006097560975609756
$ws.Range("C3").Value = 0.006097560975609756
$ws.Range("J3").Value = 0.01219512195121951
$ws.Range("P3").Value = 0.8048780487804879
$ws.Range("S3").Value = 0.1707317073170732
$ws.Range("B6").Value = 0.06451612903225806
$ws.Range("D6").Value = 0.009216589861751152
$ws.Range("F6").Value = 0.07373271889400922
$ws.Range("J6").Value = 0.2304147465437788
$ws.Range("O6").Value = 0.0184331797235023
$ws.Range("Q6").Value = 0.2211981566820277
$ws.Range("R6").Value = 0.05990783410138249
$ws.Range("S6").Value = 0.3225806451612903
$ws.Range("B7").Value = 0.091324200913242
$ws.Range("D7").Value = 0.0091324200913242
$ws.Range("E7").Value = 0.0045662100456621
$ws.Range("F7").Value = 0.0867579908675799
$ws.Range("J7").Value = 0.1187214611872146
$ws.Range("O7").Value = 0.0136986301369863
$ws.Range("Q7").Value = 0.2009132420091324
$ws.Range("R7").Value = 0.0776255707762557
$ws.Range("S7").Value = 0.3972602739726027
$ws.Range("B8").Value = 0.1014150943396226
$ws.Range("D8").Value = 0.01886792452830189
$ws.Range("F8").Value = 0.06367924528301887
$ws.Range("J8").Value = 0.09433962264150944
$ws.Range("O8").Value = 0.009433962264150943
$ws.Range("Q8").Value = 0.1910377358490566
$ws.Range("R8").Value = 0.1226415094339623
$ws.Range("S8").Value = 0.3985849056603774
$ws.Range("B9").Value = 0.09392265193370165
$ws.Range("D9").Value = 0.01657458563535912
$ws.Range("F9").Value = 0.04972375690607735
$ws.Range("J9").Value = 0.138121546961326
$ws.Range("O9").Value = 0.01104972375690608
$ws.Range("Q9").Value = 0.143646408839779
$ws.Range("R9").Value = 0.09392265193370165
$ws.Range("S9").Value = 0.4530386740331492
$ws.Range("B10").Value = 0.1053555750658472
$ws.Range("D10").Value = 0.02194907813871817
$ws.Range("E10").Value = 0.000877963125548727
$ws.Range("F10").Value = 0.07726075504828797
$ws.Range("J10").Value = 0.07550482879719052
$ws.Range("O10").Value = 0.009657594381035996
$ws.Range("Q10").Value = 0.1984196663740123
$ws.Range("R10").Value = 0.09482001755926252
$ws.Range("S10").Value = 0.4161545215100966
$ws.Range("G11").Value = 0.1329305135951662
$ws.Range("J11").Value = 0.1148036253776435
$ws.Range("K11").Value = 0.1933534743202417
$ws.Range("L11").Value = 0.5468277945619335
$ws.Range("S11").Value = 0.01208459214501511
$ws.Range("G12").Value = 0.7595628415300546
$ws.Range("J12").Value = 0.2131147540983606
$ws.Range("K12").Value = 0.00546448087431694
$ws.Range("L12").Value = 0.01092896174863388
$ws.Range("S12").Value = 0.01092896174863388
$ws.Range("F15").Value = 0.005263157894736842
$ws.Range("H15").Value = 0.1578947368421053
$ws.Range("I15").Value = 0.08947368421052632
$ws.Range("J15").Value = 0.3789473684210526
$ws.Range("K15").Value = 0.08947368421052632
$ws.Range("M15").Value = 0.02105263157894737
$ws.Range("N15").Value = 0.01578947368421053
$ws.Range("O15").Value = 0.06315789473684211
$ws.Range("S15").Value = 0.1789473684210526
$ws.Range("F16").Value = 0.01595744680851064
$ws.Range("H16").Value = 0.175531914893617
$ws.Range("I16").Value = 0.04787234042553191
$ws.Range("J16").Value = 0.4095744680851064
$ws.Range("K16").Value = 0.09574468085106383
$ws.Range("M16").Value = 0.03723404255319149
$ws.Range("O16").Value = 0.09042553191489362
$ws.Range("S16").Value = 0.1276595744680851
$ws.Range("F17").Value = 0.01913875598086124
$ws.Range("H17").Value = 0.1842105263157895
$ws.Range("I17").Value = 0.07894736842105263
$ws.Range("J17").Value = 0.437799043062201
$ws.Range("K17").Value = 0.1004784688995215
$ws.Range("M17").Value = 0.0215311004784689
$ws.Range("N17").Value = 0.002392344497607655
$ws.Range("O17").Value = 0.04066985645933014
$ws.Range("S17").Value = 0.1148325358851675
$ws.Range("F18").Value = 0.009852216748768473
$ws.Range("H18").Value = 0.1428571428571428
$ws.Range("I18").Value = 0.1379310344827586
$ws.Range("J18").Value = 0.3990147783251232
$ws.Range("K18").Value = 0.1576354679802956
$ws.Range("M18").Value = 0.03940886699507389
$ws.Range("O18").Value = 0.04926108374384237
$ws.Range("S18").Value = 0.06403940886699508
$ws.Range("F19").Value = 0.01719690455717971
$ws.Range("H19").Value = 0.2218400687876182
$ws.Range("I19").Value = 0.07996560619088564
$ws.Range("J19").Value = 0.3611349957007738
$ws.Range("K19").Value = 0.1332760103181427
$ws.Range("M19").Value = 0.02407566638005159
$ws.Range("N19").Value = 0.0008598452278589854
$ws.Range("O19").Value = 0.07222699914015478
$ws.Range("S19").Value = 0.08942390369733448
